$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data (row 16) for the "editor" table (columns L:T, Table14).
# L16 stays empty just like the other rows in this table (11-14).
# M16 = Editor Version ("0.1"). Typing "0.1" directly would be interpreted
# as a number; copy the equal text from M11 (already stored as text) so
# the new cell keeps the same text type without picking up a new style.
$ws.Range("M11").Copy()
$ws.Range("M16").PasteSpecial(-4104)
$excel.CutCopyMode = 0

$ws.Range("N16").Value = "Struct Property Editor"
$ws.Range("O16").Value = 1
$ws.Range("P16").Value = 'Ability to edit a RTTR reflected structure that is provided as a rttr::variant or rttr::instance, so that the structure can be used to automatically construct an ImGui "editor" and the values be set to the correct values '
$ws.Range("Q16").Value = "High"
$ws.Range("R16").Value = "Low"
$ws.Range("S16").Value = "Engine 0.1"
$ws.Range("T16").Value = "On Hold"

# Update the sheet view's active cell/selection to match the author's.
$ws.Range("T13").Select()
